$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.373.41'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").Value = '3.509.68'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.23'
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.44'
$ws.Range("E6").Value = '  -1.05%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +5.99%  '
$ws.Range("E10").Value = '  +0.19%  '
$ws.Range("E11").Value = '  +3.86%  '
$ws.Range("D12").Value = '4.108.65'
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("E13").Value = '  +1.18%  '
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("D15").Value = '3.509.67'
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.79'
$ws.Range("E16").Value = '  +3.35%  '
$ws.Range("D17").Value = '64.357.71'
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("E18").Value = '  +1.95%  '
$ws.Range("E19").Value = '  +3.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.58'
$ws.Range("E20").Value = '  -1.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '394.91'
$ws.Range("E21").Value = '  +2.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.577'
$ws.Range("E22").Value = '  +1.26%  '
$ws.Range("D23").Value = '3.651.43'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.74'
$ws.Range("E24").Value = '  +0.99%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("E27").Value = '  +2.14%  '
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("E29").Value = '  -1.95%  '
$ws.Range("E30").Value = '  +1.31%  '
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("E32").Value = '  -7.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.157'
$ws.Range("E33").Value = '  +6.25%  '
$ws.Range("D34").Value = '3.540.91'
$ws.Range("E34").Value = '  +0.46%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  -0.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.38'
$ws.Range("E37").Value = '  +1.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.96'
$ws.Range("E38").Value = '  +1.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.55'
$ws.Range("E39").Value = '  +0.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '167.37'
$ws.Range("E40").Value = '  +2.28%  '
$ws.Range("E41").Value = '  +0.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.812'
$ws.Range("E42").Value = '  +0.50%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.35'
$ws.Range("E44").Value = '  -2.46%  '
$ws.Range("E46").Value = '  +0.61%  '
$ws.Range("E47").Value = '  -3.03%  '
$ws.Range("E48").Value = '  +0.51%  '
$ws.Range("D49").Value = '2.381.47'
$ws.Range("E49").Value = '  -3.90%  '
$ws.Range("E50").Value = '  -1.92%  '
$ws.Range("E51").Value = '  +0.09%  '
